$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.953.46"
Set-TextValue $ws.Range("D3") "1.631.34"
Set-TextValue $ws.Range("E3") "  -0.50%  "
Set-TextValue $ws.Range("E4") "  -0.02%  "
Set-TextValue $ws.Range("D5") "211.82"
Set-TextValue $ws.Range("E6") "  -0.20%  "
Set-TextValue $ws.Range("E7") "  -0.04%  "
Set-TextValue $ws.Range("D8") "23.41"
Set-TextValue $ws.Range("E8") "  -0.92%  "
Set-TextValue $ws.Range("E9") "  -1.77%  "
Set-TextValue $ws.Range("E10") "  -0.19%  "
Set-TextValue $ws.Range("E11") "  +0.66%  "
Set-TextValue $ws.Range("D12") "1.862.92"
Set-TextValue $ws.Range("E12") "  -0.48%  "
Set-TextValue $ws.Range("D13") "1.628.18"
Set-TextValue $ws.Range("E13") "  -0.82%  "
Set-TextValue $ws.Range("E14") "  -1.23%  "
Set-TextValue $ws.Range("D15") "0.562"
Set-TextValue $ws.Range("E15") "  -2.13%  "
Set-TextValue $ws.Range("D16") "65.58"
Set-TextValue $ws.Range("E16") "  -0.86%  "
Set-TextValue $ws.Range("D17") "27.946.02"
Set-TextValue $ws.Range("E17") "  +0.17%  "
Set-TextValue $ws.Range("D18") "230.81"
Set-TextValue $ws.Range("E18") "  -0.53%  "
Set-TextValue $ws.Range("E19") "  +0.10%  "
Set-TextValue $ws.Range("D20") "7.66"
Set-TextValue $ws.Range("E20") "  +0.80%  "
Set-TextValue $ws.Range("E21") "  -0.13%  "
Set-TextValue $ws.Range("D22") "10.36"
Set-TextValue $ws.Range("E22") "  -5.13%  "
Set-TextValue $ws.Range("E24") "  -1.62%  "
Set-TextValue $ws.Range("D25") "154.99"
Set-TextValue $ws.Range("E25") "  +2.12%  "
Set-TextValue $ws.Range("E26") "  +0.27%  "
Set-TextValue $ws.Range("E27") "  -0.10%  "
Set-TextValue $ws.Range("E28") "  -0.90%  "
Set-TextValue $ws.Range("E29") "  -0.07%  "
Set-TextValue $ws.Range("E30") "  -0.58%  "
Set-TextValue $ws.Range("E31") "  -0.47%  "
Set-TextValue $ws.Range("E32") "  +1.86%  "
Set-TextValue $ws.Range("D33") "1.402.02"
Set-TextValue $ws.Range("E33") "  -0.98%  "
Set-TextValue $ws.Range("E34") "  -1.56%  "
Set-TextValue $ws.Range("E35") "  +0.32%  "
Set-TextValue $ws.Range("E37") "  +0.42%  "
Set-TextValue $ws.Range("E38") "  +2.34%  "
Set-TextValue $ws.Range("D39") "0.557"
Set-TextValue $ws.Range("E39") "  +0.36%  "
Set-TextValue $ws.Range("D40") "0.866"
Set-TextValue $ws.Range("E40") "  -2.92%  "
Set-TextValue $ws.Range("E41") "  -0.24%  "
Set-TextValue $ws.Range("E42") "  -0.09%  "
Set-TextValue $ws.Range("D43") "66.58"
Set-TextValue $ws.Range("E43") "  -0.80%  "
Set-TextValue $ws.Range("D44") "1.85"
Set-TextValue $ws.Range("E44") "  +1.41%  "
Set-TextValue $ws.Range("E46") "  -0.31%  "
Set-TextValue $ws.Range("D47") "1.773.02"
Set-TextValue $ws.Range("D48") "88.19"
Set-TextValue $ws.Range("E48") "  -0.24%  "
Set-TextValue $ws.Range("E49") "  -0.14%  "
Set-TextValue $ws.Range("D51") "7.55"
Set-TextValue $ws.Range("E51") "  -1.14%  "
